$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 35

# Column A holds a date string like "10/06/2025". Force text formatting first so
# Excel does not auto-convert it into a date serial number, then restore the
# cell's style to the default "Normal" so no stray number-format style sticks
# around on the new cell (matching the unstyled cells A2:A34 in the source file).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "10/06/2025"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1521778353556913
$ws.Cells.Item($row, 3).Value = 0.8478221646443087
